# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and restore swapped BabyDogeCoin/Aave row order (rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.073.33"
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = "'1.675.32"
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = "'210.97"
$ws.Range("E5").Value = '  -3.48%  '
$ws.Range("D6").Value = "'0.5265"
$ws.Range("E6").Value = '  -5.89%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("D8").Value = "'0.2673"
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = "'0.06305"
$ws.Range("E9").Value = '  -3.33%  '
$ws.Range("D10").Value = "'21.22"
$ws.Range("E10").Value = '  -4.29%  '
$ws.Range("D11").Value = "'0.07582"
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = "'1.675.08"
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = "'4.499"
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").Value = "'0.5667"
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").Value = "'0.000008112"
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("D16").Value = "'65.59"
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = "'26.087.21"
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D19").Value = "'4.839"
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").Value = "'10.59"
$ws.Range("E20").Value = '  -3.23%  '
$ws.Range("D21").Value = "'189.05"
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").Value = "'6.183"
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = "'148.54"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = "'0.1251"
$ws.Range("E25").Value = '  -6.20%  '
$ws.Range("D26").Value = "'7.630"
$ws.Range("E26").Value = '  -3.66%  '
$ws.Range("D27").Value = "'15.91"
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").Value = "'0.06354"
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").Value = "'1.283"
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("D31").Value = "'3.532"
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("D32").Value = "'3.522"
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("D33").Value = "'1.662"
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = "'1.009"
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("D35").Value = "'2.417"
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").Value = "'0.6034"
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = "'6.139"
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("D39").Value = "'0.01616"
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").Value = "'1.093.83"
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("D41").Value = "'0.8698"
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").Value = "'100.04"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'56.96"
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = "'0.05249"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = "'7.989"
$ws.Range("E49").Value = '  -2.69%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").Value = "'5.936"
$ws.Range("E51").Value = '  -2.58%  '
